$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.707577466964722
$ws.Range("B1").Value = 4.086183547973633
$ws.Range("C1").Value = 3.182519435882568
$ws.Range("D1").Value = 1.442652583122253
$ws.Range("E1").Value = 1.022813320159912
